# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# Map: worksheet name -> @{ row = newValue }
$updates = @{
    "展览" = @{
        3  = 19290
        17 = 1492
        22 = 8079
        23 = 990
        27 = 1259
        30 = 6099
        32 = 76
        33 = 177
        35 = 294
    }
    "全部类型" = @{
        3  = 19290
        17 = 1492
        23 = 8079
        24 = 990
        28 = 1259
        33 = 6099
        35 = 76
        36 = 177
        38 = 294
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowsMap = $updates[$sheetName]
    foreach ($row in $rowsMap.Keys) {
        $newValue = $rowsMap[$row]
        $ws.Range("F$row").Value = $newValue
    }
}
